# Raw and Clean Data from SSA for October 23-26
# Adds four new rows (147-150) of "out_vars" time series data covering
# 2020-10-23 through 2020-10-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 147-149: dates stored as real date serials (A column), formatted
# with a custom yyyy-mm-dd date format, and the data columns (B:F) use a
# slightly larger font with wrapped text - matching how this batch of
# rows was appended by the data pipeline.
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row = 147; Date = 44127; B = 880775; C = 1058102; D = 330956; E = 88312; F = 22.990548096846528 },
    @{ Row = 148; Date = 44128; B = 886800; C = 1066646; D = 335143; E = 88743; F = 22.937302661253948 },
    @{ Row = 149; Date = 44129; B = 891160; C = 1072760; D = 331758; E = 88924; F = 22.897796130885588 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $r.Date
    $dateCell.NumberFormat = "yyyy\-mm\-dd;@"

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    for ($col = 2; $col -le 6; $col++) {
        $c = $ws.Cells.Item($row, $col)
        $c.Font.Size = 12
        $c.WrapText = $true
    }

    $ws.Rows.Item($row).RowHeight = 16
}

# ---------------------------------------------------------------------
# Row 150: 2020-10-26 - the date here is stored as literal text (as the
# earlier rows that predate the date-serial switch do), with plain,
# unformatted data cells.
# ---------------------------------------------------------------------

$textDateCell = $ws.Cells.Item(150, 1)
$textDateCell.NumberFormat = "@"
$textDateCell.Value = "2020-10-26"
$textDateCell.ClearFormats()

$ws.Cells.Item(150, 2).Value = 895326
$ws.Cells.Item(150, 3).Value = 1078072
$ws.Cells.Item(150, 4).Value = 328231
$ws.Cells.Item(150, 5).Value = 89171
$ws.Cells.Item(150, 6).Value = 22.89
